{"js": "// Merge the word-by-word runs of the Title, Author and Abstract\n// paragraphs into a single run each (the runs' combined text is left\n// unchanged, only the run-splitting is collapsed).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/style,items/text\");\nawait context.sync();\n\nconst targets = {\n  \"Title\": \"Questions: The scalar product\",\n  \"Author\": \"Ritwik Anand\",\n  \"Abstract\": \"A selection of questions for the study guide on the scalar product\"\n};\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const style = para.style;\n  if (Object.prototype.hasOwnProperty.call(targets, style)) {\n    para.insertText(targets[style], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Merge the word-by-word runs of the Title, Author and Abstract\n# paragraphs into a single run each (the runs' combined text is left\n# unchanged, only the run-splitting is collapsed).\n$d = $word.ActiveDocument\n\n$targets = @{\n    \"Title\"    = \"Questions: The scalar product\"\n    \"Author\"   = \"Ritwik Anand\"\n    \"Abstract\" = \"A selection of questions for the study guide on the scalar product\"\n}\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $styleName = $p.Style.NameLocal\n    if ($targets.ContainsKey($styleName)) {\n        $newText = $targets[$styleName]\n        $r = $p.Range\n        # Replacing the paragraph's full text with itself via Find/Replace\n        # collapses the paragraph's many single-word runs into one run,\n        # matching the canonical OOXML produced by the edit.\n        $r.Find.Execute($newText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n    }\n}\n"}
